$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-dimension / iaest-measure qualifiers - several columns move from
# "dimension" to "measure" after the data was reprocessed with the newly
# curated dimensions.
$ws.Range("A2").Value = "iaest-measure:temporalidad"
$ws.Range("J2").Value = "iaest-measure:mes-nombre"
$ws.Range("L2").Value = "iaest-measure:modalidad"
$ws.Range("M2").Value = "iaest-measure:dias-duracion-contrato"
$ws.Range("O2").Value = "iaest-measure:grupo-de-tipo-de-contrato"
$ws.Range("Q2").Value = "iaest-measure:sexo"

# Row 3: dim/medida marker follows the same columns.
$ws.Range("A3").Value = "medida"
$ws.Range("J3").Value = "medida"
$ws.Range("L3").Value = "medida"
$ws.Range("M3").Value = "medida"
$ws.Range("O3").Value = "medida"
$ws.Range("Q3").Value = "medida"

# Row 4: datatype marker follows the same columns.
$ws.Range("A4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"
$ws.Range("L4").Value = "xsd:int"
$ws.Range("M4").Value = "xsd:int"
$ws.Range("O4").Value = "xsd:int"
$ws.Range("Q4").Value = "xsd:int"

# Row 5: mapping file reference - only remains for the columns that are still
# true dimensions (ccaa-nombre, ano). Fully clear the rest (value + format) so
# the cells disappear from the sheet, matching the reprocessed output.
"A5", "J5", "L5", "M5", "O5", "Q5" | ForEach-Object {
    $ws.Range($_).ClearFormats()
    $ws.Range($_).Value = $null
}
